# Web UI portfolio-data refresh: bump every row's "修改时间" (last-updated)
# stamp from 202509211549 to 202509211550 across all three portfolio
# sheets. The stamp is stored as text (it is a fixed-width yyyyMMddHHmm
# code, not a real date/time value), so it is written back with a leading
# apostrophe to keep Excel from reinterpreting it as a number.

$wb = $excel.ActiveWorkbook

$oldStamp = "202509211549"
$newStamp = "'202509211550"

# "大智投资组合": timestamp in column E, data rows 2-9
$wsDaZhi = $wb.Worksheets.Item("大智投资组合")
for ($r = 2; $r -le 9; $r++) {
    $cell = $wsDaZhi.Cells.Item($r, 5)
    if ($cell.Value() -eq $oldStamp) {
        $cell.Value = $newStamp
        # keep the cell's style untouched (General) - the leading quote
        # above only exists to stop Excel re-typing the stamp as a number
        $cell.Style = "Normal"
    }
}

# "大成投资组合": timestamp in column E, data rows 2-11
$wsDaCheng = $wb.Worksheets.Item("大成投资组合")
for ($r = 2; $r -le 11; $r++) {
    $cell = $wsDaCheng.Cells.Item($r, 5)
    if ($cell.Value() -eq $oldStamp) {
        $cell.Value = $newStamp
        $cell.Style = "Normal"
    }
}

# "我的投资组合": timestamp in column G, data rows 2-13
$wsMine = $wb.Worksheets.Item("我的投资组合")
for ($r = 2; $r -le 13; $r++) {
    $cell = $wsMine.Cells.Item($r, 7)
    if ($cell.Value() -eq $oldStamp) {
        $cell.Value = $newStamp
        $cell.Style = "Normal"
    }
}
